# "Um pouco de organização" — clean up the regex/pattern strings in column E:
# drop the superfluous backslash-escapes around literal braces ( \{ \} -> { } )
# and simplify a few of the patterns themselves. The Nome/Descrição (C/D)
# columns are untouched; only the Formato (E) column text changes, plus the
# row-15 custom height (a leftover from the old, longer wrapped text) is
# cleared back to the sheet's default.
#
# NB: cell writes below are ordered deliberately (matches the order the
# strings end up appended to the shared-string table on save).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E9").Value  = '{[1-9] [:string:] [:edicao:] [:indice:] [:paginas:] [:anexos:] [:strings:]}'
$ws.Range("E11").Value = '{[0-9]+ [0-9]+ [01] [01]}'
$ws.Range("E13").Value = '{([:subIndice:]|[:folhaIndice:])*}'
$ws.Range("E15").Value = '{[:string:] ([:subIndice:]|[:folhaIndice:])* 0}'
$ws.Range("E17").Value = '{[:string:] [:id:] 1}'
$ws.Range("E19").Value = '{[:pagina:]*}'
$ws.Range("E23").Value = '{[:anexo:]*}'
$ws.Range("E21").Value = '{[:objeto:]*}'
$ws.Range("E25").Value = '{[:string:] [:string:]}'
$ws.Range("E27").Value = '{[:string:]*}'
$ws.Range("E33").Value = '{[:alinhamento:] [:strings:] 0}'
$ws.Range("E35").Value = '{[:alinhamento:] [:string:] 1}'
$ws.Range("E37").Value = '{[:imagem:] 2}'
$ws.Range("E39").Value = '{[:alinhamento:] [:nivel:] [:string:] 3}'
$ws.Range("E41").Value = '{[:altura:] 4}'
$ws.Range("E47").Value = 'GROB [0-9]+ [0-9]+ [0-9A-F]+'
$ws.Range("E45").Value = '".*"'

# Row 15 no longer needs its explicit wrapped-text height now that the
# pattern text is shorter — restore the default row height (no explicit
# height / customHeight in the saved XML), matching the surrounding rows.
$ws.Rows.Item(15).AutoFit()
